$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the columns that hold numeric-looking / percentage-looking / date-looking strings
# so Excel keeps them as text (matching original inlineStr text cells) rather than converting to numbers/dates.
$ws.Range("B2:G51").NumberFormat = "@"

$ws.Range("B2").Value = 'BNB'
$ws.Range("C2").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("D2").Value = '302.87'
$ws.Range("E2").Value = '1.12%'
$ws.Range("F2").Value = '15-1-2023'
$ws.Range("G2").Value = '1'

$ws.Range("B3").Value = 'OKB'
$ws.Range("C3").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D3").Value = '31.95'
$ws.Range("E3").Value = '0.03%'
$ws.Range("F3").Value = '15-1-2023'
$ws.Range("G3").Value = '1'

$ws.Range("B4").Value = 'HuobiToken'
$ws.Range("C4").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D4").Value = '5.284'
$ws.Range("E4").Value = '1.74%'
$ws.Range("F4").Value = '15-1-2023'
$ws.Range("G4").Value = '1'

$ws.Range("B5").Value = 'Cronos'
$ws.Range("C5").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D5").Value = '0.07480'
$ws.Range("E5").Value = '1.17%'
$ws.Range("F5").Value = '15-1-2023'
$ws.Range("G5").Value = '1'

$ws.Range("B6").Value = 'KuCoinToken'
$ws.Range("C6").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D6").Value = '7.879'
$ws.Range("E6").Value = '2.73%'
$ws.Range("F6").Value = '15-1-2023'
$ws.Range("G6").Value = '1'

$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").Value = '3.862'
$ws.Range("E7").Value = '5.54%'
$ws.Range("F7").Value = '15-1-2023'
$ws.Range("G7").Value = '1'

$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8").Value = '1.702'
$ws.Range("E8").Value = '19.08%'
$ws.Range("F8").Value = '15-1-2023'
$ws.Range("G8").Value = '1'

$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9").Value = '0.9242'
$ws.Range("E9").Value = '1.88%'
$ws.Range("F9").Value = '15-1-2023'
$ws.Range("G9").Value = '1'

$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").Value = '0.1694'
$ws.Range("E10").Value = '2.71%'
$ws.Range("F10").Value = '15-1-2023'
$ws.Range("G10").Value = '1'

$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11").Value = '0.07636'
$ws.Range("E11").Value = '-5.79%'
$ws.Range("F11").Value = '15-1-2023'
$ws.Range("G11").Value = '1'

$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12").Value = '0.08011'
$ws.Range("E12").Value = '1.70%'
$ws.Range("F12").Value = '15-1-2023'
$ws.Range("G12").Value = '1'

$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").Value = '0.03014'
$ws.Range("E13").Value = '0.11%'
$ws.Range("F13").Value = '15-1-2023'
$ws.Range("G13").Value = '1'

$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").Value = '0.09912'
$ws.Range("E14").Value = '0.97%'
$ws.Range("F14").Value = '15-1-2023'
$ws.Range("G14").Value = '1'

$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15").Value = '0.001512'
$ws.Range("E15").Value = '-0.53%'
$ws.Range("F15").Value = '15-1-2023'
$ws.Range("G15").Value = '1'

$ws.Range("B16").Value = 'CoinExToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D16").Value = '0.04625'
$ws.Range("E16").Value = '3.32%'
$ws.Range("F16").Value = '15-1-2023'
$ws.Range("G16").Value = '1'

$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").Value = '0.006367'
$ws.Range("E17").Value = '4.22%'
$ws.Range("F17").Value = '15-1-2023'
$ws.Range("G17").Value = '1'

$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = '3.447'
$ws.Range("E18").Value = '-1.96%'
$ws.Range("F18").Value = '15-1-2023'
$ws.Range("G18").Value = '1'

$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D19").Value = '2.231'
$ws.Range("E19").Value = '-0.61%'
$ws.Range("F19").Value = '15-1-2023'
$ws.Range("G19").Value = '1'

$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Value = '0.3294'
$ws.Range("E20").Value = '0.82%'
$ws.Range("F20").Value = '15-1-2023'
$ws.Range("G20").Value = '1'

$ws.Range("B21").Value = 'ProBitToken'
$ws.Range("C21").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D21").Value = '0.1349'
$ws.Range("E21").Value = '4.15%'
$ws.Range("F21").Value = '15-1-2023'
$ws.Range("G21").Value = '1'

$ws.Range("B22").Value = 'MCDex'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D22").Value = '4.584'
$ws.Range("E22").Value = '11.55%'
$ws.Range("F22").Value = '15-1-2023'
$ws.Range("G22").Value = '1'

$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").Value = '0.1555'
$ws.Range("E23").Value = '-1.33%'
$ws.Range("F23").Value = '15-1-2023'
$ws.Range("G23").Value = '1'

$ws.Range("B24").Value = 'BitKan'
$ws.Range("C24").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D24").Value = '0.001222'
$ws.Range("E24").Value = '1.88%'
$ws.Range("F24").Value = '15-1-2023'
$ws.Range("G24").Value = '1'

$ws.Range("B25").Value = 'HotbitToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D25").Value = '0.004438'
$ws.Range("E25").Value = '-0.94%'
$ws.Range("F25").Value = '15-1-2023'
$ws.Range("G25").Value = '1'

$ws.Range("B26").Value = 'NitroEx'
$ws.Range("C26").Value = 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
$ws.Range("D26").Value = '0.0001405'
$ws.Range("E26").Value = '21.20%'
$ws.Range("F26").Value = '15-1-2023'
$ws.Range("G26").Value = '1'

$ws.Range("B27").Value = 'UpBots'
$ws.Range("C27").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D27").Value = '0.0001807'
$ws.Range("E27").Value = '7.45%'
$ws.Range("F27").Value = '15-1-2023'
$ws.Range("G27").Value = '1'

$ws.Range("B28").Value = 'Spectre.aiUtilityToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/yxQ8LoZvwJ6Ee+spectreaiutilitytoken-sxut'
$ws.Range("D28").Value = '--'
$ws.Range("E28").Value = '--%'
$ws.Range("F28").Value = '15-1-2023'
$ws.Range("G28").Value = '1'

$ws.Range("B29").Value = 'LegolasExchange'
$ws.Range("C29").Value = 'https://coinranking.com/coin/zEMEnlPs_94tc+legolasexchange-lgo'
$ws.Range("D29").Value = '--'
$ws.Range("E29").Value = '--%'
$ws.Range("F29").Value = '15-1-2023'
$ws.Range("G29").Value = '1'

$ws.Range("B30").Value = 'BitZToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/nLHHPBKDJSEee+bitztoken-bz'
$ws.Range("D30").Value = '--'
$ws.Range("E30").Value = '--%'
$ws.Range("F30").Value = '15-1-2023'
$ws.Range("G30").Value = '1'

$ws.Range("B31").Value = 'Birake'
$ws.Range("C31").Value = 'https://coinranking.com/coin/dTOfofFqKQiY5+birake-bir'
$ws.Range("D31").Value = '--'
$ws.Range("E31").Value = '--%'
$ws.Range("F31").Value = '15-1-2023'
$ws.Range("G31").Value = '1'

$ws.Range("B32").Value = 'NashExchange'
$ws.Range("C32").Value = 'https://coinranking.com/coin/9LcSTo0q-+nashexchange-nex'
$ws.Range("D32").Value = '--'
$ws.Range("E32").Value = '--%'
$ws.Range("F32").Value = '15-1-2023'
$ws.Range("G32").Value = '1'

$ws.Range("B33").Value = 'AAXToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/LNePqkIhk+aaxtoken-aab'
$ws.Range("D33").Value = '--'
$ws.Range("E33").Value = '--%'
$ws.Range("F33").Value = '15-1-2023'
$ws.Range("G33").Value = '1'

$ws.Range("B34").Value = 'CenX'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V4XJUvLQb+cenx-cenx'
$ws.Range("D34").Value = '--'
$ws.Range("E34").Value = '--%'
$ws.Range("F34").Value = '15-1-2023'
$ws.Range("G34").Value = '1'

$ws.Range("B35").Value = 'BNIXToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/n194X9uHp+bnixtoken-bnix'
$ws.Range("D35").Value = '--'
$ws.Range("E35").Value = '--%'
$ws.Range("F35").Value = '15-1-2023'
$ws.Range("G35").Value = '1'

$ws.Range("B36").Value = 'Polkally'
$ws.Range("C36").Value = 'https://coinranking.com/coin/NkDWUL8F-+polkally-kally'
$ws.Range("D36").Value = '--'
$ws.Range("E36").Value = '--%'
$ws.Range("F36").Value = '15-1-2023'
$ws.Range("G36").Value = '1'

$ws.Range("B37").Value = 'Charli3'
$ws.Range("C37").Value = 'https://coinranking.com/coin/8SgjMSqUk+charli3-c3'
$ws.Range("D37").Value = '--'
$ws.Range("E37").Value = '--%'
$ws.Range("F37").Value = '15-1-2023'
$ws.Range("G37").Value = '1'

$ws.Range("B38").Value = 'BlubitexToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Y9oImHIW5+blubitextoken-bbe'
$ws.Range("D38").Value = '--'
$ws.Range("E38").Value = '--%'
$ws.Range("F38").Value = '15-1-2023'
$ws.Range("G38").Value = '1'

$ws.Range("B39").Value = 'One'
$ws.Range("C39").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D39").Value = '0.01672'
$ws.Range("E39").Value = '8.40%'
$ws.Range("F39").Value = '15-1-2023'
$ws.Range("G39").Value = '1'

$ws.Range("B40").Value = 'IDEX'
$ws.Range("C40").Value = 'https://coinranking.com/coin/ZiRElvGxqQaf+idex-idex'
$ws.Range("D40").Value = '0.04561'
$ws.Range("E40").Value = '2.70%'
$ws.Range("F40").Value = '15-1-2023'
$ws.Range("G40").Value = '1'

$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").Value = '0.007033'
$ws.Range("E41").Value = '2.06%'
$ws.Range("F41").Value = '15-1-2023'
$ws.Range("G41").Value = '1'

$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").Value = '0.1344'
$ws.Range("E42").Value = '-0.52%'
$ws.Range("F42").Value = '15-1-2023'
$ws.Range("G42").Value = '1'

$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").Value = '0.002067'
$ws.Range("E43").Value = '-4.73%'
$ws.Range("F43").Value = '15-1-2023'
$ws.Range("G43").Value = '1'

$ws.Range("B44").Value = 'LocalTraders'
$ws.Range("C44").Value = 'https://coinranking.com/coin/E6DwMU2zXb+localtraders-lct'
$ws.Range("D44").Value = '0.01337'
$ws.Range("E44").Value = '2.67%'
$ws.Range("F44").Value = '15-1-2023'
$ws.Range("G44").Value = '1'

$ws.Range("B45").Value = 'CoinLion'
$ws.Range("C45").Value = 'https://coinranking.com/coin/sot4vgRyjNXek+coinlion-lion'
$ws.Range("D45").Value = '0.00006204'
$ws.Range("E45").Value = '7.14%'
$ws.Range("F45").Value = '15-1-2023'
$ws.Range("G45").Value = '1'

$ws.Range("B46").Value = 'BOLO'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D46").Value = '0.7073'
$ws.Range("E46").Value = '-62.38%'
$ws.Range("F46").Value = '15-1-2023'
$ws.Range("G46").Value = '1'

$ws.Range("B47").Value = 'CoinbaseStockToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D47").Value = '0.01228'
$ws.Range("E47").Value = '-4.64%'
$ws.Range("F47").Value = '15-1-2023'
$ws.Range("G47").Value = '1'

$ws.Range("B48").Value = 'DigiFinexToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/rY6dWXQL4+digifinextoken-dft'
$ws.Range("D48").Value = '--'
$ws.Range("E48").Value = '--%'
$ws.Range("F48").Value = '15-1-2023'
$ws.Range("G48").Value = '1'

$ws.Range("B49").Value = 'Bitcoin2.0'
$ws.Range("C49").Value = 'https://coinranking.com/coin/tSq1ehUma+bitcoin20-xbtc2'
$ws.Range("D49").Value = '--'
$ws.Range("E49").Value = '--%'
$ws.Range("F49").Value = '15-1-2023'
$ws.Range("G49").Value = '1'

$ws.Range("B50").Value = 'CoinField'
$ws.Range("C50").Value = 'https://coinranking.com/coin/h4GpuIkN_+coinfield-cfc'
$ws.Range("D50").Value = '--'
$ws.Range("E50").Value = '--%'
$ws.Range("F50").Value = '15-1-2023'
$ws.Range("G50").Value = '1'

$ws.Range("B51").Value = 'Coinovy'
$ws.Range("C51").Value = 'https://coinranking.com/coin/6NDu4kaME+coinovy-c2f'
$ws.Range("D51").Value = '--'
$ws.Range("E51").Value = '--%'
$ws.Range("F51").Value = '15-1-2023'
$ws.Range("G51").Value = '1'
